$wb = $excel.ActiveWorkbook

# --- Strategy & Operations sheet: append new AddDataPoint(S&O) rows ---
$ws = $wb.Worksheets.Item("Strategy & Operations")

$newRows = @(
    @("Additional Data Point"),
    @("Value"),
    @("Position Vacant"),
    @("Position Missing"),
    @("Name of the Policy"),
    @("Criticality (Y/N)"),
    @("Availability (Y/N)"),
    @("Name of the Technical Platform"),
    @("Criticality (Y/N)"),
    @("Availability (Y/N)")
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = "AddDataPoint(S&O)"
    $ws.Range("B$r").Value = $newRows[$i][0]
}

# --- Update active sheet / selection state to match authored workbook view ---
$wsGlossary = $wb.Worksheets.Item("Glossary&Definitions")
$wsGlossary.Range("J29").Select()

$ws.Activate()
$ws.Range("D30").Select()
